$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old hashcode -> new hashcode pairs (matches the sharedStrings.xml diff)
$pairs = @(
    ,@("2a8524da19a261ecdef6891100f68859", "f922ed9e447644263a1a97de707e5cb8")
    ,@("1a2aad99247432a7c8ad2c855eaeec1e", "62d0f019011e1e35afb4da08a53861dd")
    ,@("3a425473b901d99eeb2f8f05d1a7a9da", "f7945b435d376f43969ae850a7cc68cb")
    ,@("7332e19db9d80de1248db805e60f9312", "45345d734b099da46e786c83e8f28c96")
    ,@("a0cab0e46f110ea81f706b2fc5953f20", "74c498ae62afc36eaf69fb2be262b624")
    ,@("dadb7be999dbd63f806299bfafbc6261", "6b15316edc1cc092b4abac42be90bd28")
    ,@("25264021f32130c246ff1dcdeec483d0", "a971ea9eb8c3823f3586968e3793190b")
    ,@("b2958ca0a2f48c38ed413b0942283382", "d9cbdf45e33118bc240620a3976be092")
    ,@("6988a7d8079cd0690a573f7b9e2adedb", "d7a63990157d9dcc566c9c52d107a4bf")
    ,@("9ff250cc2296e8b04e2e9c55eb7b492a", "ae42a0af0e2092a422639ad4d71db265")
    ,@("289d9c7f686850f0271f00b042591a5a", "618db607106c4c865cbafcf8156b579a")
    ,@("d0a510c33c0ac6bb6a7521f08fec4070", "b2c2d7b0c6e1e482e2baebfaa3e80238")
    ,@("82760c335d1800fd1aeb50687d6f826e", "811e4b110a2cffba77fce045c7017d73")
    ,@("9a8cc75de1629534c3eaece5b8c32057", "67e8de9238b1d980854c534789e8446c")
    ,@("3f0a589ba5292d038af5d7e15f995d2b", "869c621bbced2dd1e9009bcaac137d49")
    ,@("181895aa68478a8ce5e37e3a6123fdf6", "beba7bce29c4068483cd10898052ff4a")
    ,@("c16252edd9bbad81bece7e1e437aeca5", "7d3192fea74a6be1ead9e53c83c35f0f")
    ,@("930e9bd628ccd09c643cd2b4a4b8cfad", "0841f66eec1f7caf51680bed6f5054c6")
    ,@("76fe75e6b689c434da60d249ba6765bf", "e1e4b714dddf2e3deb6075c4d94ffcf9")
    ,@("090ce60a84e4df080ad7c313bf00d29a", "7c7e26fef28b133513b0e1d817db11ed")
    ,@("b3c0471f6ab03fe79ed3515cd46b22cc", "3bb24bf20af84bd73d4fd48e30da03f3")
    ,@("3573f972709eca56275fd504bb286c75", "0f2b68cdf56bae47118f70f03e78d2f5")
    ,@("e3d6f2571a6e47a237de56acc60583d0", "3962d32114f3fb69ae6f12f86a119019")
    ,@("f329d36fa47d84734dfb9b3626f9d4ca", "46abcc7d85f2732d753478da077c6dad")
    ,@("c2ff6a83c1beba8689e2d6eaa3eb06e1", "320c9d5b1e38d46bf285d4beb72f820c")
    ,@("6872b106d46507f66af37d33523f76f9", "caed40e30b8d326c9ee29159f49801d9")
    ,@("a43aad2a42277be6fc85233bafe81f21", "94c8a699ba72fa2ba49483e62eaeeb5b")
    ,@("2ba2af195a7150411e9edbf214040e44", "6dae6fa19d878e3e786208dc34f13627")
    ,@("db79560a07b943a028661bf9ac58f8cf", "0500c3294f2fe90971052abfee60871b")
    ,@("16b63d480f3d50d78a869c19ab998727", "2ede366eee4394e48ea0925f9464345c")
    ,@("7f37c26eae181fa0ad2e97b5864751b2", "654c1ba0472b17af82efd250300ae113")
    ,@("4fc5fa4b3dd3ce2d2f863a4ac7f1255b", "87f7d8c8d5f14748512c9245c79f6ea6")
    ,@("536052429b70078e1e780ef554fbc516", "e992428de39ad6cc52cb72f089587295")
    ,@("32cabfb6d54c47197f02bfa132f2bceb", "c73244e4d02da93b2f5418460dd36c9d")
    ,@("fb3404a2ee3af1938e8f92d2e045b730", "d174fa8fbca0c777f41402c2571309ad")
    ,@("14cb8d34718c47516b19ad2970bcf17c", "c3305368066951b035b3eec49bbfc9ce")
    ,@("93049bfcc2ff1ccbc37fcd3a7fe75f92", "1f9b18a75e7137204200fd2e581624f2")
    ,@("73dcb4033cf74069e3da205ee99500a5", "bebe597650251d7dc4b5abfc624cebb2")
    ,@("7b32c2e2138ad20d6de90800ca768f42", "0a647b4a3f32e50bca26867df944df5e")
    ,@("1240d1925d5bb6781d888325f1408e49", "d05f60cb7fe7ed68b218c83ac767a514")
    ,@("18959c8132fbe58132b63e2ed262ede7", "828dfcdbe017b46b27ba6a91372baea2")
    ,@("683ad9d5a62eedccab952d06bed5a4f7", "369163dccc3c430a954a07963037cfd1")
    ,@("c23d1d2e9e89bd032e026d27dfcc8827", "55ee70e9919cf8142a528225a340560d")
    ,@("97010d418992034607b9ffb8ac4a8020", "e8dfad8ff97156163b1440cb8b6475c6")
)

$notFound = @()
foreach ($pair in $pairs) {
    $oldVal = $pair[0]
    $newVal = $pair[1]
    $found = $ws.Cells.Find($oldVal)
    if ($found) {
        $found.Value = $newVal
    } else {
        $notFound += $oldVal
    }
}

Write-Host "Updated $($pairs.Count - $notFound.Count) of $($pairs.Count) hashcodes."
if ($notFound.Count -gt 0) {
    Write-Host "Not found:"
    foreach ($nf in $notFound) {
        Write-Host "  $nf"
    }
}
